$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the helper text in D1: a space was inserted after "P," so the
# comma-separated list now reads "A,C,D,H,J,P, R,S,V" instead of
# "A,C,D,H,J,P,R,S,V".
$ws.Range("D1").Value = "A,C,D,H,J,P, R,S,V"

# Move/leave the active selection on D1 (was D6 before the edit).
$ws.Range("D1").Select()
